{"js": "// Blockquotes (the \"Block Text\" style) should be single-spaced.\n// Give the \"Block Text\" paragraph style explicit single line spacing\n// (renders as line=\"240\" lineRule=\"auto\" in the OOXML) plus 12pt (240 twips)\n// space-after, matching Word's standard \"single spacing\" paragraph preset.\nconst style = context.document.getStyles().getByNameOrNullObject(\"Block Text\");\nawait context.sync();\n\nif (!style.isNullObject) {\n  style.paragraphFormat.lineSpacing = 12; // points -> single line spacing\n  style.paragraphFormat.spaceAfter = 12;  // points -> 240 twips\n  await context.sync();\n}\n", "ps1": "# Blockquotes (the \"Block Text\" style) should be single-spaced.\n# Give the \"Block Text\" paragraph style explicit single line spacing\n# (line=\"240\" lineRule=\"auto\") plus 12pt (240 twips) space-after, matching\n# Word's standard \"single spacing\" paragraph spacing preset.\n$d = $word.ActiveDocument\n$s = $d.Styles(\"Block Text\")\n$s.ParagraphFormat.LineSpacingRule = 0   # wdLineSpaceSingle\n$s.ParagraphFormat.SpaceAfter = 12       # points (240 twips)\n"}
